$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = "'29.030.04"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.04%  "
$ws.Cells.Item(3, 4).Value = "'1.833.13"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "
$ws.Cells.Item(4, 4).Value = "'0.9986"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "'242.32"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.35%  "
$ws.Cells.Item(6, 4).Value = "'0.6273"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -4.25%  "
$ws.Cells.Item(7, 4).Value = "'0.9995"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "
$ws.Cells.Item(8, 4).Value = "'0.07620"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.98%  "
$ws.Cells.Item(9, 4).Value = "'0.2928"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.08%  "
$ws.Cells.Item(10, 4).Value = "'22.54"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.87%  "
$ws.Cells.Item(11, 4).Value = "'0.07712"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.45%  "
$ws.Cells.Item(12, 4).Value = "'1.839.96"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.37%  "
$ws.Cells.Item(13, 4).Value = "'4.956"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.30%  "
$ws.Cells.Item(14, 4).Value = "'0.6654"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.11%  "
$ws.Cells.Item(15, 4).Value = "'0.00001014"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +16.98%  "
$ws.Cells.Item(16, 4).Value = "'82.77"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.59%  "
$ws.Cells.Item(17, 4).Value = "'6.047"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.05%  "
$ws.Cells.Item(18, 4).Value = "'29.038.10"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.06%  "
$ws.Cells.Item(19, 4).Value = "'226.54"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.07%  "
$ws.Cells.Item(20, 5).Value = "  -0.60%  "
$ws.Cells.Item(21, 5).Value = "  -0.10%  "
$ws.Cells.Item(22, 4).Value = "'7.180"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.81%  "
$ws.Cells.Item(23, 4).Value = "'0.9996"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.10%  "
$ws.Cells.Item(24, 4).Value = "'158.28"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.51%  "
$ws.Cells.Item(25, 4).Value = "'8.499"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.13%  "
$ws.Cells.Item(26, 5).Value = "  -0.50%  "
$ws.Cells.Item(27, 4).Value = "'17.91"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.08%  "
$ws.Cells.Item(28, 4).Value = "'1.489"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.07%  "
$ws.Cells.Item(29, 4).Value = "'4.111"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.14%  "
$ws.Cells.Item(30, 4).Value = "'4.018"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.24%  "
$ws.Cells.Item(31, 4).Value = "'1.191"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.69%  "
$ws.Cells.Item(32, 4).Value = "'0.05225"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.68%  "
$ws.Cells.Item(33, 5).Value = "  +0.20%  "
$ws.Cells.Item(34, 4).Value = "'0.7356"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.73%  "
$ws.Cells.Item(35, 4).Value = "'1.139"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.43%  "
$ws.Cells.Item(36, 4).Value = "'2.707"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +2.16%  "
$ws.Cells.Item(37, 4).Value = "'1.240.76"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.40%  "
$ws.Cells.Item(38, 5).Value = "  -0.01%  "
$ws.Cells.Item(39, 4).Value = "'0.01786"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.09%  "
$ws.Cells.Item(40, 4).Value = "'6.348"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.10%  "
$ws.Cells.Item(41, 4).Value = "'0.8950"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.56%  "
$ws.Cells.Item(42, 4).Value = "'0.9997"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.00%  "
$ws.Cells.Item(43, 4).Value = "'101.78"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.25%  "
$ws.Cells.Item(44, 4).Value = "'1.981.97"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.35%  "
$ws.Cells.Item(45, 5).Value = "  -2.08%  "
$ws.Cells.Item(46, 4).Value = "'64.20"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.60%  "
$ws.Cells.Item(48, 4).Value = "'0.4038"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.39%  "
$ws.Cells.Item(49, 4).Value = "'8.933"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.20%  "
$ws.Cells.Item(50, 4).Value = "'1.642"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -5.60%  "
$ws.Cells.Item(51, 4).Value = "'0.05747"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.53%  "
